# Update selected cells in the "Anual" worksheet to reflect the latest
# seasonally-adjusted data revision from MV (Actualizacion desde MV -datos-).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    # Row 14
    $ws.Range("B14").Value = 107370
    $ws.Range("D14").Value = 65617
    $ws.Range("H14").Value = 13777
    $ws.Range("L14").Value = 41467
    $ws.Range("M14").Value = 35518
    $ws.Range("O14").Value = 20060
    $ws.Range("P14").Value = 12901
    $ws.Range("X14").Value = 113737
    # Row 15
    $ws.Range("B15").Value = 100404
    $ws.Range("D15").Value = 64983
    $ws.Range("L15").Value = 39716
    $ws.Range("M15").Value = 34475
    $ws.Range("O15").Value = 20283
    $ws.Range("P15").Value = 11795
    $ws.Range("X15").Value = 111987
    # Row 16
    $ws.Range("B16").Value = 113993
    $ws.Range("D16").Value = 71887
    $ws.Range("L16").Value = 40650
    $ws.Range("M16").Value = 34377
    $ws.Range("O16").Value = 20478
    $ws.Range("X16").Value = 118555
    # Row 17
    $ws.Range("B17").Value = 124691
    $ws.Range("F17").Value = 32348
    $ws.Range("L17").Value = 42843
    $ws.Range("M17").Value = 35776
    $ws.Range("O17").Value = 20133
    $ws.Range("X17").Value = 125692
    # Row 18
    $ws.Range("B18").Value = 133886
    $ws.Range("C18").Value = 99029
    $ws.Range("D18").Value = 82570
    $ws.Range("F18").Value = 34176
    $ws.Range("L18").Value = 42923
    $ws.Range("M18").Value = 36450
    $ws.Range("O18").Value = 20575
    $ws.Range("X18").Value = 132514
    # Row 19
    $ws.Range("B19").Value = 138723
    $ws.Range("L19").Value = 44424
    $ws.Range("M19").Value = 38001
    $ws.Range("P19").Value = 13570
    $ws.Range("S19").Value = 39185
    $ws.Range("V19").Value = 33796
    $ws.Range("X19").Value = 137929
    # Row 20
    $ws.Range("B20").Value = 138023
    $ws.Range("I20").Value = 32567
    $ws.Range("J20").Value = 20886
    $ws.Range("K20").Value = 11679
    $ws.Range("L20").Value = 44551
    $ws.Range("M20").Value = 38616
    $ws.Range("O20").Value = 22042
    $ws.Range("X20").Value = 140306
    # Row 21
    $ws.Range("B21").Value = 141466
    $ws.Range("C21").Value = 109024
    $ws.Range("D21").Value = 90575
    $ws.Range("J21").Value = 21499
    $ws.Range("L21").Value = 43744
    $ws.Range("M21").Value = 38088
    $ws.Range("R21").Value = 41757
    $ws.Range("S21").Value = 36398
    $ws.Range("V21").Value = 31339
    $ws.Range("X21").Value = 143475
    # Row 22
    $ws.Range("B22").Value = 143950
    $ws.Range("C22").Value = 112802
    $ws.Range("K22").Value = 10770
    $ws.Range("L22").Value = 43942
    $ws.Range("M22").Value = 38220
    $ws.Range("O22").Value = 21191
    $ws.Range("P22").Value = 13821
    $ws.Range("X22").Value = 145801
    # Row 23
    $ws.Range("B23").Value = 148461
    $ws.Range("F23").Value = 40008
    $ws.Range("I23").Value = 31034
    $ws.Range("J23").Value = 20108
    $ws.Range("K23").Value = 10919
    $ws.Range("L23").Value = 43259
    $ws.Range("M23").Value = 37471
    $ws.Range("O23").Value = 20513
    $ws.Range("P23").Value = 13891
    $ws.Range("R23").Value = 44083
    $ws.Range("X23").Value = 147852
    # Row 24
    $ws.Range("B24").Value = 155105
    $ws.Range("I24").Value = 32602
    $ws.Range("J24").Value = 20701
    $ws.Range("K24").Value = 11882
    $ws.Range("L24").Value = 45614
    $ws.Range("M24").Value = 39828
    $ws.Range("O24").Value = 21577
    $ws.Range("P24").Value = 14603
    $ws.Range("R24").Value = 47651
    $ws.Range("S24").Value = 41922
    $ws.Range("V24").Value = 36245
    $ws.Range("X24").Value = 153419
    # Row 25
    $ws.Range("B25").Value = 156526
    $ws.Range("C25").Value = 122331
    $ws.Range("D25").Value = 100952
    $ws.Range("F25").Value = 41208
    $ws.Range("H25").Value = 21378
    $ws.Range("I25").Value = 34056
    $ws.Range("J25").Value = 21855
    $ws.Range("K25").Value = 12183
    $ws.Range("L25").Value = 44457
    $ws.Range("P25").Value = 14164
    $ws.Range("R25").Value = 46524
    $ws.Range("S25").Value = 40880
    $ws.Range("V25").Value = 35058
    $ws.Range("X25").Value = 154766
    # Row 26
    $ws.Range("B26").Value = 141995
    $ws.Range("C26").Value = 113899
    $ws.Range("D26").Value = 93245
    $ws.Range("F26").Value = 40151
    $ws.Range("H26").Value = 20589
    $ws.Range("I26").Value = 30081
    $ws.Range("J26").Value = 19346
    $ws.Range("L26").Value = 43959
    $ws.Range("M26").Value = 40051
    $ws.Range("P26").Value = 14787
    $ws.Range("R26").Value = 40578
    $ws.Range("S26").Value = 36688
    $ws.Range("V26").Value = 31161
    $ws.Range("X26").Value = 145498
